$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A3").Value = "Cost"
$ws.Range("B3").Value = 2

$ws.Range("C3").Select()
